$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Coinranking prices/volumes are plain text in this sheet (often numeric-
    # looking, e.g. "582.65"). Excel's normal type-inference would silently
    # turn those into real numbers, so force the Text category first for any
    # value that looks numeric; values that already contain things like a
    # second "." (e.g. "71.161.93") are never number-like and don't need it.
    $cell = $ws.Range($range)
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

# --- Rows 2..41: in-place price / volume updates (no row/col/name shuffling) ---

Set-TextValue "D2" "71.161.93"
Set-TextValue "E2" "  +2.61%  "

Set-TextValue "D3" "3.692.13"
Set-TextValue "E3" "  +8.06%  "

Set-TextValue "E4" "  +0.02%  "

Set-TextValue "D5" "582.65"
Set-TextValue "E5" "  +0.15%  "

Set-TextValue "D6" "177.41"
Set-TextValue "E6" "  +0.70%  "

Set-TextValue "D7" "3.682.83"
Set-TextValue "E7" "  +7.98%  "

Set-TextValue "E8" "  +4.15%  "

Set-TextValue "E9" "  -0.03%  "

Set-TextValue "E10" "  +0.78%  "

Set-TextValue "D11" "6.84"
Set-TextValue "E11" "  +26.64%  "

Set-TextValue "E12" "  +4.75%  "

Set-TextValue "D13" "48.90"
Set-TextValue "E13" "  +0.69%  "

Set-TextValue "E14" "  +2.29%  "

Set-TextValue "D15" "4.286.39"
Set-TextValue "E15" "  +8.05%  "

Set-TextValue "D16" "676.99"
Set-TextValue "E16" "  -1.75%  "

Set-TextValue "E17" "  +4.66%  "

Set-TextValue "D18" "3.682.37"
Set-TextValue "E18" "  +7.74%  "

Set-TextValue "D19" "71.327.00"
Set-TextValue "E19" "  +2.72%  "

Set-TextValue "E20" "  +1.16%  "

Set-TextValue "D21" "17.93"
Set-TextValue "E21" "  +2.00%  "

Set-TextValue "D22" "11.58"
Set-TextValue "E22" "  +2.30%  "

Set-TextValue "D23" "0.940"
Set-TextValue "E23" "  +5.48%  "

Set-TextValue "D24" "17.37"
Set-TextValue "E24" "  +3.14%  "

Set-TextValue "D25" "101.87"
Set-TextValue "E25" "  +0.79%  "

Set-TextValue "D26" "3.97"
Set-TextValue "E26" "  +2.03%  "

Set-TextValue "E27" "  +5.94%  "

Set-TextValue "E28" "  +7.73%  "

Set-TextValue "D29" "35.07"
Set-TextValue "E29" "  +5.28%  "

Set-TextValue "E30" "  +5.44%  "

Set-TextValue "E31" "  -1.45%  "

Set-TextValue "E32" "  +6.56%  "

Set-TextValue "E33" "  +11.20%  "

Set-TextValue "D34" "586.95"
Set-TextValue "E34" "  +1.87%  "

Set-TextValue "D35" "11.18"
Set-TextValue "E35" "  +1.93%  "

Set-TextValue "E36" "  +5.60%  "

Set-TextValue "D37" "58.96"
Set-TextValue "E37" "  +1.28%  "

Set-TextValue "E38" "  +0.12%  "

Set-TextValue "D39" "3.676.50"
Set-TextValue "E39" "  +3.27%  "

Set-TextValue "E40" "  -0.12%  "

Set-TextValue "E41" "  +5.02%  "

# --- Rows 42..51: coin list reshuffle (TheGraph moves up, ApeXProtocol drops out,
#     LidoDAOToken is newly appended) plus matching price/volume updates ---

Set-TextValue "B42" "TheGraph"
Set-TextValue "C42" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D42" "0.349"
Set-TextValue "E42" "  +5.58%  "

Set-TextValue "B43" "InjectiveProtocol"
Set-TextValue "C43" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D43" "35.41"
Set-TextValue "E43" "  +2.00%  "

Set-TextValue "B44" "PEPE"
Set-TextValue "C44" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D44" "0.0₃0763"
Set-TextValue "E44" "  +5.37%  "

Set-TextValue "B45" "Stacks"
Set-TextValue "C45" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "3.42"
Set-TextValue "E45" "  +5.51%  "

Set-TextValue "B46" "Fetch.AI"
Set-TextValue "C46" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D46" "2.75"
Set-TextValue "E46" "  +4.00%  "

Set-TextValue "B47" "VeChain"
Set-TextValue "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0455"
Set-TextValue "E47" "  +9.53%  "

Set-TextValue "E48" "  +9.15%  "

Set-TextValue "B49" "Stellar"
Set-TextValue "C49" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.133"
Set-TextValue "E49" "  +3.56%  "

Set-TextValue "B50" "Monero"
Set-TextValue "C50" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "135.56"
Set-TextValue "E50" "  +2.36%  "

Set-TextValue "B51" "LidoDAOToken"
Set-TextValue "C51" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D51" "2.96"
Set-TextValue "E51" "  +9.85%  "
